$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diversity")
$wsTaxa = $wb.Worksheets.Item("Ichthyoplankton Taxa")

$regions = @("MAB","GB","GOM")
$fvals = @(1.64818021721719,2.4311043939187802,2.3106921453177298,2.6576968241933399,1.8890217633458499,2.3083940791697302,2.6762277337174001,2.58663250369499,2.2392431396125998,1.92800880629282,2.50979148748229,2.0724496112749802,2.47699369010643,2.5136662597141499,2.4115712965910601,2.1840716245113301,2.1452559353519298,2.3238299006724001,1.83412114895976,2.20353217928117,2.1644707676561499,0.953180798659644,1.7186384848318901,2.0099489651570002,1.86712329301817,2.2190713181937798,1.58523211889477,1.83781990475627,2.2736607734929302,0.80185571950397305,0.78163387698259601,2.03189819099219,2.0948832890626798,1.7587995301359101,2.2952788471172099,1.9974103699685599,1.5535218818854999,2.0575666441182698,2.07657227799915,1.9426549688914201,2.0695966627581401,1.93673233663968,2.2773581889587899,1.79262114622353,2.1162296671332701,2.3616977925468001,1.4620638538459001,1.7424912417337901,1.99949717344967,2.1749223880627202,2.1250134746903999,1.8382406840057699,1.77017086236796,2.1912922562933699,2.0516064959836502,1.7892832331287101,2.1969654976724602,1.7249140146911399,1.8871028619130401,0.37327561959218702,1.5338774971665099,2.0902600976740402,1.3276940604644001)

$row = 2
for ($b = 0; $b -lt 3; $b++) {
    $region = $regions[$b]
    for ($y = 0; $y -lt 21; $y++) {
        $year = 1999 + $y
        $idx = ($b * 21) + $y
        $ws.Cells.Item($row, 1).Value = $year
        $ws.Cells.Item($row, 2).Value = "Ich_Shannon-Wiener_Diversity_index"
        $ws.Cells.Item($row, 3).Value = "Unitless"
        $ws.Cells.Item($row, 4).Value = $region
        $ws.Cells.Item($row, 5).Value = "NEFSC EcoMon Data"
        $ws.Cells.Item($row, 6).Value = $fvals[$idx]
        $row = $row + 1
    }
}

$ws.Activate()
$ws.Range("F44").Select()
$wsTaxa.Select()
$ws.Select()
